# Auto-generated Excel COM-interop script
# Applies league-base update for Australia ALeague dated 31-03-2024
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 137 is a brand new row - copy cell styles from row 136 first
$ws.Range("A136").Copy($ws.Range("A137"))
$ws.Range("E136").Copy($ws.Range("E137"))

# Row 104
$ws.Cells.Item(104, 1).Value = 102
$ws.Cells.Item(104, 2).Value = 7127370
$ws.Cells.Item(104, 3).Value = "Australia ALeague"
$ws.Cells.Item(104, 4).Value = "Australia ALeague"
$ws.Cells.Item(104, 5).Value = 45340.125
$ws.Cells.Item(104, 6).Value = "Macarthur FC"
$ws.Cells.Item(104, 7).Value = "Wellington Phoenix"
$ws.Cells.Item(104, 8).Value = 1
$ws.Cells.Item(104, 9).Value = 2
$ws.Cells.Item(104, 10).Value = "A"
$ws.Cells.Item(104, 11).Value = 2.4
$ws.Cells.Item(104, 12).Value = 3.75
$ws.Cells.Item(104, 13).Value = 2.625
$ws.Cells.Item(104, 14).Value = 2.375
$ws.Cells.Item(104, 15).Value = 3.8
$ws.Cells.Item(104, 16).Value = 2.75
$ws.Cells.Item(104, 17).Value = 0
$ws.Cells.Item(104, 18).Value = 1.8
$ws.Cells.Item(104, 19).Value = 2.05
$ws.Cells.Item(104, 20).Value = 3
$ws.Cells.Item(104, 21).Value = 1.9
$ws.Cells.Item(104, 22).Value = 1.95
$ws.Cells.Item(104, 23).Value = -1
$ws.Cells.Item(104, 24).Value = -1
$ws.Cells.Item(104, 25).Value = 1.75
$ws.Cells.Item(104, 26).Value = -1
$ws.Cells.Item(104, 27).Value = 1.05
$ws.Cells.Item(104, 28).Value = 0
$ws.Cells.Item(104, 29).Value = -0

# Row 105
$ws.Cells.Item(105, 1).Value = 103
$ws.Cells.Item(105, 2).Value = 7127374
$ws.Cells.Item(105, 3).Value = "Australia ALeague"
$ws.Cells.Item(105, 4).Value = "Australia ALeague"
$ws.Cells.Item(105, 5).Value = 45340.125
$ws.Cells.Item(105, 6).Value = "Central Coast Mariners"
$ws.Cells.Item(105, 7).Value = "Western Sydney Wanderers"
$ws.Cells.Item(105, 8).Value = 1
$ws.Cells.Item(105, 9).Value = 0
$ws.Cells.Item(105, 10).Value = "H"
$ws.Cells.Item(105, 11).Value = 1.909
$ws.Cells.Item(105, 12).Value = 3.75
$ws.Cells.Item(105, 13).Value = 3.6
$ws.Cells.Item(105, 14).Value = 2.15
$ws.Cells.Item(105, 15).Value = 3.6
$ws.Cells.Item(105, 16).Value = 3.25
$ws.Cells.Item(105, 17).Value = -0.25
$ws.Cells.Item(105, 18).Value = 1.86
$ws.Cells.Item(105, 19).Value = 2.04
$ws.Cells.Item(105, 20).Value = 2.75
$ws.Cells.Item(105, 21).Value = 1.975
$ws.Cells.Item(105, 22).Value = 1.875
$ws.Cells.Item(105, 23).Value = 1.15
$ws.Cells.Item(105, 24).Value = -1
$ws.Cells.Item(105, 25).Value = -1
$ws.Cells.Item(105, 26).Value = 0.8600000000000001
$ws.Cells.Item(105, 27).Value = -1
$ws.Cells.Item(105, 28).Value = -1
$ws.Cells.Item(105, 29).Value = 0.875

# Row 112
$ws.Cells.Item(112, 1).Value = 110
$ws.Cells.Item(112, 2).Value = 7127379
$ws.Cells.Item(112, 3).Value = "Australia ALeague"
$ws.Cells.Item(112, 4).Value = "Australia ALeague"
$ws.Cells.Item(112, 5).Value = 45347.125
$ws.Cells.Item(112, 6).Value = "Melbourne Victory"
$ws.Cells.Item(112, 7).Value = "Central Coast Mariners"
$ws.Cells.Item(112, 8).Value = 0
$ws.Cells.Item(112, 9).Value = 1
$ws.Cells.Item(112, 10).Value = "A"
$ws.Cells.Item(112, 11).Value = 1.95
$ws.Cells.Item(112, 12).Value = 3.6
$ws.Cells.Item(112, 13).Value = 3.8
$ws.Cells.Item(112, 14).Value = 1.909
$ws.Cells.Item(112, 15).Value = 3.6
$ws.Cells.Item(112, 16).Value = 4
$ws.Cells.Item(112, 17).Value = -0.5
$ws.Cells.Item(112, 18).Value = 1.9
$ws.Cells.Item(112, 19).Value = 1.95
$ws.Cells.Item(112, 20).Value = 2.75
$ws.Cells.Item(112, 21).Value = 1.925
$ws.Cells.Item(112, 22).Value = 1.925
$ws.Cells.Item(112, 23).Value = -1
$ws.Cells.Item(112, 24).Value = -1
$ws.Cells.Item(112, 25).Value = 3
$ws.Cells.Item(112, 26).Value = -1
$ws.Cells.Item(112, 27).Value = 0.95
$ws.Cells.Item(112, 28).Value = -1
$ws.Cells.Item(112, 29).Value = 0.925

# Row 113
$ws.Cells.Item(113, 1).Value = 111
$ws.Cells.Item(113, 2).Value = 7127376
$ws.Cells.Item(113, 3).Value = "Australia ALeague"
$ws.Cells.Item(113, 4).Value = "Australia ALeague"
$ws.Cells.Item(113, 5).Value = 45347.125
$ws.Cells.Item(113, 6).Value = "Newcastle Jets"
$ws.Cells.Item(113, 7).Value = "Macarthur FC"
$ws.Cells.Item(113, 8).Value = 2
$ws.Cells.Item(113, 9).Value = 2
$ws.Cells.Item(113, 10).Value = "D"
$ws.Cells.Item(113, 11).Value = 1.95
$ws.Cells.Item(113, 12).Value = 4
$ws.Cells.Item(113, 13).Value = 3.4
$ws.Cells.Item(113, 14).Value = 1.909
$ws.Cells.Item(113, 15).Value = 4.2
$ws.Cells.Item(113, 16).Value = 3.6
$ws.Cells.Item(113, 17).Value = -0.5
$ws.Cells.Item(113, 18).Value = 1.89
$ws.Cells.Item(113, 19).Value = 2.01
$ws.Cells.Item(113, 20).Value = 3.5
$ws.Cells.Item(113, 21).Value = 1.95
$ws.Cells.Item(113, 22).Value = 1.9
$ws.Cells.Item(113, 23).Value = -1
$ws.Cells.Item(113, 24).Value = 3.2
$ws.Cells.Item(113, 25).Value = -1
$ws.Cells.Item(113, 26).Value = -1
$ws.Cells.Item(113, 27).Value = 1.01
$ws.Cells.Item(113, 28).Value = 0.95
$ws.Cells.Item(113, 29).Value = -1

# Row 133
$ws.Cells.Item(133, 1).Value = 131
$ws.Cells.Item(133, 2).Value = 7126793
$ws.Cells.Item(133, 3).Value = "Australia ALeague"
$ws.Cells.Item(133, 4).Value = "Australia ALeague"
$ws.Cells.Item(133, 5).Value = 45381.14583333334
$ws.Cells.Item(133, 6).Value = "Melbourne City"
$ws.Cells.Item(133, 7).Value = "Newcastle Jets"
$ws.Cells.Item(133, 8).Value = 0
$ws.Cells.Item(133, 9).Value = 0
$ws.Cells.Item(133, 10).Value = "D"
$ws.Cells.Item(133, 11).Value = 1.571
$ws.Cells.Item(133, 12).Value = 4.333
$ws.Cells.Item(133, 13).Value = 5
$ws.Cells.Item(133, 14).Value = 1.4
$ws.Cells.Item(133, 15).Value = 5.25
$ws.Cells.Item(133, 16).Value = 6.5
$ws.Cells.Item(133, 17).Value = -1.25
$ws.Cells.Item(133, 18).Value = 1.825
$ws.Cells.Item(133, 19).Value = 2.025
$ws.Cells.Item(133, 20).Value = 3.5
$ws.Cells.Item(133, 21).Value = 1.975
$ws.Cells.Item(133, 22).Value = 1.875
$ws.Cells.Item(133, 23).Value = -1
$ws.Cells.Item(133, 24).Value = 4.25
$ws.Cells.Item(133, 25).Value = -1
$ws.Cells.Item(133, 26).Value = -1
$ws.Cells.Item(133, 27).Value = 1.025
$ws.Cells.Item(133, 28).Value = -1
$ws.Cells.Item(133, 29).Value = 0.875

# Row 134
$ws.Cells.Item(134, 1).Value = 132
$ws.Cells.Item(134, 2).Value = 7127396
$ws.Cells.Item(134, 3).Value = "Australia ALeague"
$ws.Cells.Item(134, 4).Value = "Australia ALeague"
$ws.Cells.Item(134, 5).Value = 45381.23958333334
$ws.Cells.Item(134, 6).Value = "Sydney FC"
$ws.Cells.Item(134, 7).Value = "Central Coast Mariners"
$ws.Cells.Item(134, 8).Value = 2
$ws.Cells.Item(134, 9).Value = 0
$ws.Cells.Item(134, 10).Value = "H"
$ws.Cells.Item(134, 11).Value = 2.15
$ws.Cells.Item(134, 12).Value = 3.6
$ws.Cells.Item(134, 13).Value = 3.1
$ws.Cells.Item(134, 14).Value = 2.3
$ws.Cells.Item(134, 15).Value = 3.6
$ws.Cells.Item(134, 16).Value = 2.9
$ws.Cells.Item(134, 17).Value = -0.25
$ws.Cells.Item(134, 18).Value = 2.025
$ws.Cells.Item(134, 19).Value = 1.825
$ws.Cells.Item(134, 20).Value = 2.75
$ws.Cells.Item(134, 21).Value = 1.85
$ws.Cells.Item(134, 22).Value = 2
$ws.Cells.Item(134, 23).Value = 1.3
$ws.Cells.Item(134, 24).Value = -1
$ws.Cells.Item(134, 25).Value = -1
$ws.Cells.Item(134, 26).Value = 1.025
$ws.Cells.Item(134, 27).Value = -1
$ws.Cells.Item(134, 28).Value = -1
$ws.Cells.Item(134, 29).Value = 1

# Row 135
$ws.Cells.Item(135, 1).Value = 133
$ws.Cells.Item(135, 2).Value = 7127398
$ws.Cells.Item(135, 3).Value = "Australia ALeague"
$ws.Cells.Item(135, 4).Value = "Australia ALeague"
$ws.Cells.Item(135, 5).Value = 45383.04166666666
$ws.Cells.Item(135, 6).Value = "Macarthur FC"
$ws.Cells.Item(135, 7).Value = "Western Sydney Wanderers"
$ws.Cells.Item(135, 11).Value = 2.5
$ws.Cells.Item(135, 12).Value = 3.5
$ws.Cells.Item(135, 13).Value = 2.625
$ws.Cells.Item(135, 14).Value = 2.625
$ws.Cells.Item(135, 15).Value = 3.8
$ws.Cells.Item(135, 16).Value = 2.45
$ws.Cells.Item(135, 17).Value = 0
$ws.Cells.Item(135, 18).Value = 2.02
$ws.Cells.Item(135, 19).Value = 1.88
$ws.Cells.Item(135, 20).Value = 3.25
$ws.Cells.Item(135, 21).Value = 1.925
$ws.Cells.Item(135, 22).Value = 1.925
$ws.Cells.Item(135, 23).Value = 0
$ws.Cells.Item(135, 24).Value = 0
$ws.Cells.Item(135, 25).Value = 0
$ws.Cells.Item(135, 26).Value = 0
$ws.Cells.Item(135, 27).Value = 0

# Row 136
$ws.Cells.Item(136, 1).Value = 134
$ws.Cells.Item(136, 2).Value = 7898681
$ws.Cells.Item(136, 3).Value = "Australia ALeague"
$ws.Cells.Item(136, 4).Value = "Australia ALeague"
$ws.Cells.Item(136, 5).Value = 45384.20833333334
$ws.Cells.Item(136, 6).Value = "Central Coast Mariners"
$ws.Cells.Item(136, 7).Value = "Melbourne City"
$ws.Cells.Item(136, 11).Value = 2.1
$ws.Cells.Item(136, 12).Value = 4
$ws.Cells.Item(136, 13).Value = 3
$ws.Cells.Item(136, 14).Value = 2.15
$ws.Cells.Item(136, 15).Value = 3.8
$ws.Cells.Item(136, 16).Value = 2.9
$ws.Cells.Item(136, 17).Value = -0.25
$ws.Cells.Item(136, 18).Value = 1.95
$ws.Cells.Item(136, 19).Value = 1.95
$ws.Cells.Item(136, 20).Value = 3
$ws.Cells.Item(136, 21).Value = 2.05
$ws.Cells.Item(136, 22).Value = 1.8
$ws.Cells.Item(136, 23).Value = 0
$ws.Cells.Item(136, 24).Value = 0
$ws.Cells.Item(136, 25).Value = 0
$ws.Cells.Item(136, 26).Value = 0
$ws.Cells.Item(136, 27).Value = 0

# Row 137
$ws.Cells.Item(137, 1).Value = 135
$ws.Cells.Item(137, 2).Value = 7661947
$ws.Cells.Item(137, 3).Value = "Australia ALeague"
$ws.Cells.Item(137, 4).Value = "Australia ALeague"
$ws.Cells.Item(137, 5).Value = 45385.32291666666
$ws.Cells.Item(137, 6).Value = "Perth Glory"
$ws.Cells.Item(137, 7).Value = "Sydney FC"
$ws.Cells.Item(137, 11).Value = 3.1
$ws.Cells.Item(137, 12).Value = 3.6
$ws.Cells.Item(137, 13).Value = 2.2
$ws.Cells.Item(137, 14).Value = 3.8
$ws.Cells.Item(137, 15).Value = 4.2
$ws.Cells.Item(137, 16).Value = 1.8
$ws.Cells.Item(137, 17).Value = 0.5
$ws.Cells.Item(137, 18).Value = 2.04
$ws.Cells.Item(137, 19).Value = 1.86
$ws.Cells.Item(137, 20).Value = 3.25
$ws.Cells.Item(137, 21).Value = 1.875
$ws.Cells.Item(137, 22).Value = 1.975
$ws.Cells.Item(137, 23).Value = 0
$ws.Cells.Item(137, 24).Value = 0
$ws.Cells.Item(137, 25).Value = 0
$ws.Cells.Item(137, 26).Value = 0
$ws.Cells.Item(137, 27).Value = 0

